# Updates the "cryptos" price list with freshly scraped price / 1h-volume
# figures (and reorders two rows whose underlying ranking swapped).
#
# Column D ("Price") and column E ("Volume(1h)") are stored as *text* in
# this sheet (e.g. "39.358.81", "  +1.54%  "), not as numbers - the sheet
# uses dotted thousands-separators and no real numeric semantics are
# needed. When a plain-looking numeric string (e.g. "228.69") is assigned
# through .Value, Excel's normal type inference would turn it into a
# genuine number, so those assignments are prefixed with a leading
# apostrophe to force a text value, exactly like typing it into Excel by
# hand. Non-numeric-looking text (the multi-dot prices, the percentage
# strings with their padding, names and URLs) is assigned as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells to write: cell reference -> new text value.
$updates = [ordered]@{
    'D2'  = '39.374.23'
    'E2'  = '  +1.57%  '

    'D3'  = '2.159.54'
    'E3'  = '  +3.12%  '

    'E4'  = '  +0.00%  '

    'D5'  = "'228.69"
    'E5'  = '  -0.40%  '

    'E6'  = '  +0.89%  '

    'D7'  = "'64.19"
    'E7'  = '  +5.02%  '

    'E9'  = '  +2.77%  '

    'D10' = "'0.0859"
    'E10' = '  +1.74%  '

    'E11' = '  -0.12%  '

    'D12' = "'15.93"
    'E12' = '  +3.56%  '

    'D13' = '2.479.43'
    'E13' = '  +0.53%  '

    'D14' = "'22.24"
    'E14' = '  +1.09%  '

    'E15' = '  +1.18%  '

    'E16' = '  +1.13%  '

    'D17' = '2.158.77'
    'E17' = '  +2.89%  '

    'D18' = '39.321.43'
    'E18' = '  +1.55%  '

    'D19' = "'71.91"
    'E19' = '  -0.13%  '

    'E20' = '  +0.96%  '

    'D21' = '0.0₃0851'
    'E21' = '  +1.12%  '

    'D22' = "'230.99"
    'E22' = '  +1.39%  '

    'E23' = '  +0.00%  '

    'E24' = '  +6.27%  '

    'E25' = '  +0.83%  '

    # Rows 26/27 swapped places (Monero now ranks above Cosmos).
    'B26' = 'Monero'
    'C26' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D26' = "'171.95"
    'E26' = '  +0.36%  '

    'B27' = 'Cosmos'
    'C27' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D27' = "'9.49"
    'E27' = '  -0.53%  '

    'D28' = "'0.138"
    'E28' = '  -0.35%  '

    'D29' = "'19.94"
    'E29' = '  +3.35%  '

    'E30' = '  -1.78%  '

    'E31' = '  +9.62%  '

    'E32' = '  +1.09%  '

    'E33' = '  +2.36%  '

    'E34' = '  +1.13%  '

    'D35' = "'7.08"
    'E35' = '  +7.33%  '

    'E36' = '  +1.42%  '

    'D37' = "'2.42"
    'E37' = '  +1.05%  '

    'D38' = "'3.58"
    'E38' = '  -0.50%  '

    'E39' = '  -0.22%  '

    'D40' = "'103.98"
    'E40' = '  +2.92%  '

    'D41' = "'0.0229"
    'E41' = '  +0.01%  '

    'D42' = "'17.80"
    'E42' = '  -1.93%  '

    'D43' = '1.538.09'
    'E43' = '  +0.34%  '

    'E44' = '  +3.81%  '

    'D45' = "'4.31"
    'E45' = '  +4.71%  '

    'E46' = '  +0.44%  '

    'D47' = "'0.0925"
    'E47' = '  +1.00%  '

    'D48' = "'1.10"
    'E48' = '  +5.63%  '

    'D49' = "'7.77"
    'E49' = '  +1.24%  '

    'D50' = '2.362.66'
    'E50' = '  +3.14%  '

    'E51' = '  -0.38%  '
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
